$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Heures TD") to hold the new
# "Groupes CM" information. This shifts the former F:I ("Heures TD",
# "Groupes TD", "Heures TP", "Groupes TP") into G:J.
$ws.Columns("F:F").Insert()

# Header for the newly inserted column.
$ws.Range("F1").Value = "Groupes CM"

# Populate "Groupes CM" values for each data row (2-15). Most rows use a
# single CM group; rows 4 and 12 (ALG0331 / ALG0431) use three.
$groupesCM = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 3
    13 = 1
    14 = 1
    15 = 1
}

foreach ($row in $groupesCM.Keys) {
    $ws.Cells.Item($row, 6).Value = $groupesCM[$row]
}

# Refresh the sort definition (used for the table sortState) so that it
# covers the new column layout (A2:J33 instead of A2:I33) while keeping
# the original two-level sort (by Code EC, then Code Apogée).
$sortObj = $ws.Sort
$sortFields = $sortObj.SortFields
$sortFields.Clear()
$sortFields.Add($ws.Range("B2:B33"))
$sortFields.Add($ws.Range("A2:A33"))
$sortObj.SetRange($ws.Range("A2:J33"))
$sortObj.Apply()

# Update the view: scroll so column B is the left-most visible column and
# select cell E14.
$win = $excel.ActiveWindow
$win.Zoom = 160
$win.ScrollColumn = 2
$ws.Range("E14").Select() | Out-Null
